$d = $word.ActiveDocument

# Helper: replace a whole paragraph's text with a single clean run,
# dropping any w:proofErr markers that Word's spell/grammar checker had
# stashed in the paragraph. We delete the paragraph range (which removes
# the paragraph mark and any proofErr siblings living directly in the
# paragraph), insert a fresh paragraph mark, then fill in the new text.
function Set-CleanParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $rng.Delete()
    $rng.InsertParagraphBefore()
    $rng.InsertBefore($newText)
}

# Process paragraphs from the bottom of the document upward so earlier
# edits don't shift the indices of paragraphs we haven't processed yet.
Set-CleanParagraphText 19 "Php artisan"
Set-CleanParagraphText 16 "Php artisan db:seed"
Set-CleanParagraphText 15 "Use the model name, import and define inside Model page class,"
Set-CleanParagraphText 12 "Define type using fake()"
Set-CleanParagraphText 10 "Php artisan migrate:fresh --seed"
Set-CleanParagraphText 8 "Php artisan migrate"
Set-CleanParagraphText 5 "Php artisan make model Name -a"
Set-CleanParagraphText 3 "Php artisan key:generate"

# "Command --h" -> split into two runs: keep "Command --h" as-is and
# append a second run containing "elp" (so "Command --h" + "elp" reads
# as "Command --help" but stays as two <w:r> elements). A plain
# InsertAfter would get silently coalesced into the preceding run since
# they'd share identical formatting, so we drop a temporary bookmark at
# the split point first to force the two runs to stay distinct, then
# remove the bookmark.
$p = $d.Paragraphs($d.Paragraphs.Count)
$rng = $p.Range
$rng.End = $rng.End - 1
$rng.Collapse(0)
$d.Bookmarks.Add("zzSplitPoint", $rng)
$insertRng = $rng.Duplicate
$insertRng.InsertAfter("elp")
$d.Bookmarks("zzSplitPoint").Delete()
